$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.857.28"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").Value = "3.495.72"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.79"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.87"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.02%  "
$ws.Range("D7").Value = "3.496.18"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.09"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "4.099.98"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.70"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "3.494.13"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").Value = "63.942.60"
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.98"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.40"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.64"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.35"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.578"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("D24").Value = "3.641.75"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.99"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000112"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.49"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.96%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.25"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("D33").Value = "3.500.40"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.86"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.144"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.31"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.99"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.92"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0809"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.09"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.811"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.89"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.38"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("E47").Value = "  -6.57%  "
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.89"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.436.32"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.897"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.55%  "
